$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that should remain plain text even when it looks numeric
# (mirrors typing into a cell already formatted as Text in Excel - keeps the
# original default "General" style afterwards so no stray number format sticks).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "69.918.96"
$ws.Range("E2").Value = "  +1.32%  "

$ws.Range("D3").Value = "3.518.34"
$ws.Range("E3").Value = "  -0.01%  "

Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  +0.06%  "

Set-TextValue $ws.Range("D5") "605.91"
$ws.Range("E5").Value = "  +4.40%  "

Set-TextValue $ws.Range("D6") "170.87"
$ws.Range("E6").Value = "  -2.00%  "

$ws.Range("D8").Value = "3.514.89"
$ws.Range("E8").Value = "  +0.10%  "

$ws.Range("E9").Value = "  -0.13%  "

Set-TextValue $ws.Range("D10") "0.201"
$ws.Range("E10").Value = "  +6.16%  "

$ws.Range("E11").Value = "  -0.71%  "

Set-TextValue $ws.Range("D12") "0.582"
$ws.Range("E12").Value = "  -2.90%  "

Set-TextValue $ws.Range("D13") "47.30"
$ws.Range("E13").Value = "  +0.06%  "

$ws.Range("E14").Value = "  +0.80%  "

$ws.Range("D15").Value = "4.080.98"
$ws.Range("E15").Value = "  +0.02%  "

Set-TextValue $ws.Range("D16") "8.35"
$ws.Range("E16").Value = "  -5.35%  "

Set-TextValue $ws.Range("D17") "616.53"
$ws.Range("E17").Value = "  -8.31%  "

$ws.Range("D18").Value = "3.514.55"
$ws.Range("E18").Value = "  +0.06%  "

$ws.Range("D19").Value = "69.886.83"
$ws.Range("E19").Value = "  +1.39%  "

$ws.Range("E20").Value = "  -1.90%  "

Set-TextValue $ws.Range("D21") "17.29"
$ws.Range("E21").Value = "  -1.52%  "

Set-TextValue $ws.Range("D22") "10.14"
$ws.Range("E22").Value = "  -9.84%  "

Set-TextValue $ws.Range("D23") "0.883"
$ws.Range("E23").Value = "  -2.70%  "

$ws.Range("E24").Value = "  -3.07%  "

Set-TextValue $ws.Range("D25") "96.05"
$ws.Range("E25").Value = "  -2.49%  "

Set-TextValue $ws.Range("D26") "3.87"
$ws.Range("E26").Value = "  +0.22%  "

$ws.Range("E27").Value = "  +0.07%  "

Set-TextValue $ws.Range("D28") "2.61"
$ws.Range("E28").Value = "  -2.12%  "

$ws.Range("E29").Value = "  -2.13%  "

Set-TextValue $ws.Range("D30") "33.20"
$ws.Range("E30").Value = "  +0.65%  "

$ws.Range("E31").Value = "  -3.74%  "

Set-TextValue $ws.Range("D32") "3.09"
$ws.Range("E32").Value = "  -4.15%  "

$ws.Range("E33").Value = "  -1.89%  "

Set-TextValue $ws.Range("D34") "6.96"
$ws.Range("E34").Value = "  -5.01%  "

Set-TextValue $ws.Range("D35") "562.42"
$ws.Range("E35").Value = "  -2.73%  "

Set-TextValue $ws.Range("D36") "10.77"
$ws.Range("E36").Value = "  -1.68%  "

Set-TextValue $ws.Range("D37") "3.53"
$ws.Range("E37").Value = "  -1.83%  "

Set-TextValue $ws.Range("D38") "57.08"
$ws.Range("E38").Value = "  -0.13%  "

$ws.Range("E39").Value = "  -3.63%  "

$ws.Range("E40").Value = "  +0.03%  "

Set-TextValue $ws.Range("D41") "0.141"
$ws.Range("E41").Value = "  +3.66%  "

Set-TextValue $ws.Range("D42") "0.0451"
$ws.Range("E42").Value = "  +2.59%  "

$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue $ws.Range("D43") "0.327"
$ws.Range("E43").Value = "  -3.56%  "

$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "3.333.14"
$ws.Range("E44").Value = "  -2.54%  "

Set-TextValue $ws.Range("D45") "33.04"
$ws.Range("E45").Value = "  -1.34%  "

$ws.Range("D46").Value = "0.0₃0705"
$ws.Range("E46").Value = "  -0.27%  "

$ws.Range("E47").Value = "  +0.35%  "

Set-TextValue $ws.Range("D48") "2.62"
$ws.Range("E48").Value = "  +0.70%  "

$ws.Range("E49").Value = "  -3.48%  "

Set-TextValue $ws.Range("D50") "136.46"
$ws.Range("E50").Value = "  +3.75%  "

$ws.Range("E51").Value = "  +7.31%  "
